$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The account-statement table (rows 16-21, columns C/D/E = Doc Number / Worker Name / Period)
# is refreshed: previous workers/periods are removed and new ones (for
# ALEJANDRO BERMUDEZ FERNANDEZ and CLARA INES GUZMAN MARTINEZ, periods
# 2406/2407/2408) are written in, per "Elimina EC anteriores y se agregan
# nuevos, se modifica base de datos".

# ALEJANDRO BERMUDEZ FERNANDEZ - periods 2408, 2407, 2406
$ws.Range("C16").Value = "71264684"
$ws.Range("D16").Value = "ALEJANDRO BERMUDEZ FERNANDEZ"
$ws.Range("E16").Value = "2408"

$ws.Range("C17").Value = "71264684"
$ws.Range("D17").Value = "ALEJANDRO BERMUDEZ FERNANDEZ"
$ws.Range("E17").Value = "2407"

$ws.Range("C18").Value = "71264684"
$ws.Range("D18").Value = "ALEJANDRO BERMUDEZ FERNANDEZ"
$ws.Range("E18").Value = "2406"

# CLARA INES GUZMAN MARTINEZ - periods 2408, 2407, 2406
$ws.Range("C19").Value = "45761241"
$ws.Range("D19").Value = "CLARA INES GUZMAN MARTINEZ"
$ws.Range("E19").Value = "2408"

$ws.Range("C20").Value = "45761241"
$ws.Range("D20").Value = "CLARA INES GUZMAN MARTINEZ"
$ws.Range("E20").Value = "2407"

$ws.Range("C21").Value = "45761241"
$ws.Range("D21").Value = "CLARA INES GUZMAN MARTINEZ"
$ws.Range("E21").Value = "2406"
